# Daily attendance processing - reorder "Recorded By" (column G) entries.
# For every row on the active sheet whose column G value is a
# comma-separated list of recorders, rotate the list left by one
# position (the first entry moves to the end), except for the specific
# combination "admin@admin.com, System" which is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "" -and $val -like "*, *" -and $val -ne "admin@admin.com, System") {
        $parts = $val -split ", "
        if ($parts.Length -gt 1) {
            $rotated = @()
            for ($i = 1; $i -lt $parts.Length; $i++) {
                $rotated += $parts[$i]
            }
            $rotated += $parts[0]
            $cell.Value2 = $rotated -join ", "
        }
    }
}
